# Feature improvements for renters() class
#
# Updates the "Damage State" values across the renters/owners/forsale_stock/
# forrent_stock sheets, tweaks the insurance figures for two owners, fixes
# an occupancy typo on the for-rent stock sheet, and moves the active
# sheet/selection back to the "renters" sheet.

$wb = $excel.ActiveWorkbook

# --- renters ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("renters")
$ws1.Range("M2").Value = "Extensive"
$ws1.Range("M3").Value = "Complete"
$ws1.Range("M4").Value = "Complete"
$ws1.Range("M5").Value = "Complete"

# --- owners -----------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("owners")
$ws2.Range("D3").Value = 0.85
$ws2.Range("D5").Value = 0.85
$ws2.Range("M3").Value = "None"
$ws2.Range("M4").Value = "None"
$ws2.Range("M5").Value = "None"

# --- forsale_stock -----------------------------------------------------------
$ws3 = $wb.Worksheets.Item("forsale_stock")
$ws3.Range("I2").Value = "None"
$ws3.Range("I3").Value = "None"
$ws3.Range("I4").Value = "None"
$ws3.Range("I5").Value = "None"

# --- forrent_stock -----------------------------------------------------------
$ws4 = $wb.Worksheets.Item("forrent_stock")
$ws4.Range("B2").Value = "Mobile Home"
$ws4.Range("I2").Value = "Complete"
$ws4.Range("I3").Value = "Complete"
$ws4.Range("I4").Value = "Complete"
$ws4.Range("I5").Value = "Complete"

# --- sheet view selections ---------------------------------------------------
$ws2.Range("H12").Select()
$ws3.Range("J13").Select()
$ws4.Range("I2:I5").Select()

$ws7 = $wb.Worksheets.Item("Sheet1")
$ws7.Range("J27").Select()

# Make "renters" the active sheet/tab with the new active cell selection.
$ws1.Activate()
$ws1.Range("M2").Select()
